$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ex")

# Update the "Validtill" dates
$ws.Range("F7").Value = 44466
$ws.Range("F12").Value = 44464

# Move the active selection from I9 to I12
[void]$ws.Range("I12").Select()
